$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4369
$ws.Range("I86").Value = 2121.75
$ws.Range("J86").Value = 5186.1816
$ws.Range("K86").Value = 2121.75
$ws.Range("L86").Value = 5186.1816
$ws.Range("M86").Value = -998.75
$ws.Range("N86").Value = -7432.1816

$ws.Range("H88").Value = 1931.4445
$ws.Range("I88").Value = 995.6667
$ws.Range("J88").Value = 2399.3333
$ws.Range("K88").Value = 995.6667
$ws.Range("L88").Value = 2399.3333
$ws.Range("M88").Value = -589.6667
$ws.Range("N88").Value = -3211.3333

$ws.Range("H89").Value = 4369
$ws.Range("I89").Value = 2121.75
$ws.Range("J89").Value = 5186.1816
$ws.Range("K89").Value = 10608.75
$ws.Range("L89").Value = 25930.908
$ws.Range("M89").Value = -4992.75
$ws.Range("N89").Value = -37162.908

$ws.Range("H91").Value = 1931.4445
$ws.Range("I91").Value = 995.6667
$ws.Range("J91").Value = 2399.3333
$ws.Range("K91").Value = 995.6667
$ws.Range("L91").Value = 2399.3333
$ws.Range("M91").Value = 408.3333
$ws.Range("N91").Value = -5207.3333

$ws.Range("H98").Value = 798.5454999999999
$ws.Range("I98").Value = 798.5454999999999
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 798.5454999999999
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 699.4545000000001

$ws.Range("H122").Value = 798.5454999999999
$ws.Range("I122").Value = 798.5454999999999
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2395.6365
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 54.36350000000039

$ws.Range("H132").Value = 1682.5
$ws.Range("I132").Value = 1765.9412
$ws.Range("J132").Value = 1398.8
$ws.Range("K132").Value = 5297.8236
$ws.Range("L132").Value = 4196.4
$ws.Range("M132").Value = -2767.8236
$ws.Range("N132").Value = -9256.4

$ws.Range("H141").Value = 3959.1667
$ws.Range("I141").Value = 3601
$ws.Range("J141").Value = 5750
$ws.Range("K141").Value = 10803
$ws.Range("L141").Value = 17250
$ws.Range("M141").Value = -5623
$ws.Range("N141").Value = -27610

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1406.0984
$ws.Range("I61").Value = 1403.017
$ws.Range("J61").Value = 1497
$ws.Range("K61").Value = 1403.017
$ws.Range("L61").Value = 1497
$ws.Range("M61").Value = -1191.017
$ws.Range("N61").Value = -1921

$ws.Range("H97").Value = 866.5
$ws.Range("I97").Value = 849.8
$ws.Range("J97").Value = 950
$ws.Range("K97").Value = 849.8
$ws.Range("L97").Value = 950
$ws.Range("M97").Value = -353.8
$ws.Range("N97").Value = -1942

$ws.Range("H122").Value = 528678.4
$ws.Range("I122").Value = 590546.7
$ws.Range("J122").Value = 2797.5
$ws.Range("K122").Value = 1771640.1
$ws.Range("L122").Value = 8392.5
$ws.Range("M122").Value = -1769190.1
$ws.Range("N122").Value = -13292.5

$ws.Range("H136").Value = 1406.0984
$ws.Range("I136").Value = 1403.017
$ws.Range("J136").Value = 1497
$ws.Range("K136").Value = 4209.051
$ws.Range("L136").Value = 4491
$ws.Range("M136").Value = -1659.051
$ws.Range("N136").Value = -9591

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 590.06665
$ws.Range("I22").Value = 600.0714
$ws.Range("J22").Value = 450
$ws.Range("K22").Value = 600.0714
$ws.Range("L22").Value = 450
$ws.Range("M22").Value = -427.0714
$ws.Range("N22").Value = -796

$ws.Range("H86").Value = 5268
$ws.Range("I86").Value = 4498
$ws.Range("J86").Value = 5653
$ws.Range("K86").Value = 4498
$ws.Range("L86").Value = 5653
$ws.Range("M86").Value = -3375
$ws.Range("N86").Value = -7899

$ws.Range("H89").Value = 5268
$ws.Range("I89").Value = 4498
$ws.Range("J89").Value = 5653
$ws.Range("K89").Value = 22490
$ws.Range("L89").Value = 28265
$ws.Range("M89").Value = -16874
$ws.Range("N89").Value = -39497

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3340.647
$ws.Range("I58").Value = 1862.6
$ws.Range("J58").Value = 5452.143
$ws.Range("K58").Value = 1862.6
$ws.Range("L58").Value = 5452.143
$ws.Range("M58").Value = -1659.6
$ws.Range("N58").Value = -5858.143

$ws.Range("H122").Value = 8424.286
$ws.Range("I122").Value = 8454.888999999999
$ws.Range("J122").Value = 8369.200000000001
$ws.Range("K122").Value = 25364.667
$ws.Range("L122").Value = 25107.6
$ws.Range("M122").Value = -22914.667
$ws.Range("N122").Value = -30007.6

$ws.Range("H132").Value = 2503.2334
$ws.Range("I132").Value = 2030.7084
$ws.Range("J132").Value = 4393.3335
$ws.Range("K132").Value = 6092.1252
$ws.Range("L132").Value = 13180.0005
$ws.Range("M132").Value = -3562.1252
$ws.Range("N132").Value = -18240.0005

$ws.Range("H136").Value = 3340.647
$ws.Range("I136").Value = 1862.6
$ws.Range("J136").Value = 5452.143
$ws.Range("K136").Value = 5587.799999999999
$ws.Range("L136").Value = 16356.429
$ws.Range("M136").Value = -3037.799999999999
$ws.Range("N136").Value = -21456.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 654.2308
$ws.Range("I5").Value = 395.875
$ws.Range("J5").Value = 1067.6
$ws.Range("K5").Value = 1187.625
$ws.Range("L5").Value = 3202.8
$ws.Range("M5").Value = -1075.625
$ws.Range("N5").Value = -3426.8

$ws.Range("H7").Value = 7142933.5
$ws.Range("I7").Value = 10000057
$ws.Range("J7").Value = 125
$ws.Range("K7").Value = 30000171
$ws.Range("L7").Value = 375
$ws.Range("M7").Value = -30000059
$ws.Range("N7").Value = -599

$ws.Range("H17").Value = 151
$ws.Range("I17").Value = 161.6
$ws.Range("J17").Value = 137.75
$ws.Range("K17").Value = 484.8
$ws.Range("L17").Value = 413.25
$ws.Range("M17").Value = -315.8
$ws.Range("N17").Value = -751.25

$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = ""
$ws.Range("N52").Value = 0

$ws.Range("H86").Value = 339.5
$ws.Range("I86").Value = 285
$ws.Range("J86").Value = 394
$ws.Range("K86").Value = 855
$ws.Range("L86").Value = 1182
$ws.Range("M86").Value = 331
$ws.Range("N86").Value = -3554

$ws.Range("H89").Value = 339.5
$ws.Range("I89").Value = 285
$ws.Range("J89").Value = 394
$ws.Range("K89").Value = 2565
$ws.Range("L89").Value = 3546
$ws.Range("M89").Value = 3363
$ws.Range("N89").Value = -15402

$ws.Range("H129").Value = 2079.1
$ws.Range("I129").Value = 1828
$ws.Range("J129").Value = 2665
$ws.Range("K129").Value = 5484
$ws.Range("L129").Value = 7995
$ws.Range("M129").Value = -484
$ws.Range("N129").Value = -17995

$ws.Range("H131").Value = 3241.5881
$ws.Range("I131").Value = 3479.1428
$ws.Range("J131").Value = 3075.3
$ws.Range("K131").Value = 10437.4284
$ws.Range("L131").Value = 9225.900000000001
$ws.Range("M131").Value = -5397.428400000001
$ws.Range("N131").Value = -19305.9

$ws.Range("H135").Value = 654.2308
$ws.Range("I135").Value = 395.875
$ws.Range("J135").Value = 1067.6
$ws.Range("K135").Value = 3562.875
$ws.Range("L135").Value = 9608.4
$ws.Range("M135").Value = -1027.875
$ws.Range("N135").Value = -14678.4

$ws.Range("H140").Value = 2476.6924
$ws.Range("I140").Value = 2476.6924
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 7430.0772
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -2250.0772

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6758.0713
$ws.Range("I70").Value = 6161
$ws.Range("J70").Value = 7089.778
$ws.Range("K70").Value = 6161
$ws.Range("L70").Value = 7089.778
$ws.Range("M70").Value = -5891
$ws.Range("N70").Value = -7629.778

$ws.Range("H73").Value = 6758.0713
$ws.Range("I73").Value = 6161
$ws.Range("J73").Value = 7089.778
$ws.Range("K73").Value = 6161
$ws.Range("L73").Value = 7089.778
$ws.Range("M73").Value = -5225
$ws.Range("N73").Value = -8961.778

$ws.Range("H80").Value = 2713.4167
$ws.Range("I80").Value = 3095.2856
$ws.Range("J80").Value = 2178.8
$ws.Range("K80").Value = 3095.2856
$ws.Range("L80").Value = 2178.8
$ws.Range("M80").Value = -2097.2856
$ws.Range("N80").Value = -4174.8

$ws.Range("H83").Value = 2713.4167
$ws.Range("I83").Value = 3095.2856
$ws.Range("J83").Value = 2178.8
$ws.Range("K83").Value = 15476.428
$ws.Range("L83").Value = 10894
$ws.Range("M83").Value = -10484.428
$ws.Range("N83").Value = -20878

$ws.Range("H132").Value = 3263.12
$ws.Range("I132").Value = 2835
$ws.Range("J132").Value = 4975.6
$ws.Range("K132").Value = 8505
$ws.Range("L132").Value = 14926.8
$ws.Range("M132").Value = -5975
$ws.Range("N132").Value = -19986.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2248.25
$ws.Range("I7").Value = 2248.25
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2248.25
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = -2136.25

$ws.Range("H22").Value = 966.6667
$ws.Range("I22").Value = 959.6
$ws.Range("J22").Value = 1002
$ws.Range("K22").Value = 959.6
$ws.Range("L22").Value = 1002
$ws.Range("M22").Value = -664.6
$ws.Range("N22").Value = -1592

$ws.Range("H27").Value = 966.6667
$ws.Range("I27").Value = 959.6
$ws.Range("J27").Value = 1002
$ws.Range("K27").Value = 959.6
$ws.Range("L27").Value = 1002
$ws.Range("M27").Value = -852.6
$ws.Range("N27").Value = -1216

$ws.Range("H61").Value = 2127.9285
$ws.Range("I61").Value = 1526.2727
$ws.Range("J61").Value = 4334
$ws.Range("K61").Value = 1526.2727
$ws.Range("L61").Value = 4334
$ws.Range("M61").Value = -1324.2727
$ws.Range("N61").Value = -4738

$ws.Range("H93").Value = 350
$ws.Range("I93").Value = 300
$ws.Range("J93").Value = 500
$ws.Range("K93").Value = 300
$ws.Range("L93").Value = 500
$ws.Range("M93").Value = 948
$ws.Range("N93").Value = -2996

$ws.Range("H113").Value = 2127.9285
$ws.Range("I113").Value = 1526.2727
$ws.Range("J113").Value = 4334
$ws.Range("K113").Value = 1526.2727
$ws.Range("L113").Value = 4334
$ws.Range("M113").Value = 643.7273
$ws.Range("N113").Value = -8674

$ws.Range("H122").Value = 304
$ws.Range("I122").Value = 304
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 912
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 1538

$ws.Range("H126").Value = 2248.25
$ws.Range("I126").Value = 2248.25
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6744.75
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = -4274.75

$ws.Range("H132").Value = 4784.8535
$ws.Range("I132").Value = 4344.5
$ws.Range("J132").Value = 5733.3076
$ws.Range("K132").Value = 13033.5
$ws.Range("L132").Value = 17199.9228
$ws.Range("M132").Value = -10503.5
$ws.Range("N132").Value = -22259.9228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1641.0834
$ws.Range("I81").Value = 1799.8334
$ws.Range("J81").Value = 1482.3334
$ws.Range("K81").Value = 3599.6668
$ws.Range("L81").Value = 2964.6668
$ws.Range("M81").Value = -2538.6668
$ws.Range("N81").Value = -5086.6668

$ws.Range("H84").Value = 1641.0834
$ws.Range("I84").Value = 1799.8334
$ws.Range("J84").Value = 1482.3334
$ws.Range("K84").Value = 17998.334
$ws.Range("L84").Value = 14823.334
$ws.Range("M84").Value = -12694.334
$ws.Range("N84").Value = -25431.334

$ws.Range("H96").Value = 1584.5
$ws.Range("I96").Value = 1527.5
$ws.Range("J96").Value = 1670
$ws.Range("K96").Value = 1527.5
$ws.Range("L96").Value = 1670
$ws.Range("M96").Value = -154.5
$ws.Range("N96").Value = -4416

$ws.Range("H122").Value = 1373.4
$ws.Range("I122").Value = 1373.4
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4120.200000000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1670.200000000001

$ws.Range("H136").Value = 54277.156
$ws.Range("I136").Value = 1022.4286
$ws.Range("J136").Value = 203390.4
$ws.Range("K136").Value = 3067.2858
$ws.Range("L136").Value = 610171.2
$ws.Range("M136").Value = -517.2857999999997
$ws.Range("N136").Value = -615271.2
